# Generate Report for Handoff
# Adds a new localization-status row (for the new file
# ffff4be74221-375a-4846-b937-d78ae62f0c2f.md) to the Overview, zh-cn and
# de-de sheets, and refreshes the existing row's renamed source file
# (cdd2fe4f-... -> 5be07d11-...) and timestamps.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8091b82469223b4d9966e72d2e3d2abeada4da69/e2e/"
$hyperlinkColor = 15570276  # OLE BGR for RGB FF6495ED (matches workbook's custom HyperLink style)

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows("2:2").Copy()
$wsOverview.Rows("3:3").Insert()

$wsOverview.Range("A2").Value = "5be07d11-a570-42c3-9722-9a26b9ee2579.md"
$wsOverview.Range("G2").Value = "2016-11-03 19:48:14"

$wsOverview.Range("A3").Value = "ffff4be74221-375a-4846-b937-d78ae62f0c2f.md"
$wsOverview.Range("B3").Value = "e2e\ffff4be74221-375a-4846-b937-d78ae62f0c2f.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-11-03 19:48:14"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), ($repoBase + "5be07d11-a570-42c3-9722-9a26b9ee2579.md"), "", "", "e2e\5be07d11-a570-42c3-9722-9a26b9ee2579.md")
$wsOverview.Range("B2").Font.Color = $hyperlinkColor

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($repoBase + "ffff4be74221-375a-4846-b937-d78ae62f0c2f.md"), "", "", "e2e\ffff4be74221-375a-4846-b937-d78ae62f0c2f.md")
$wsOverview.Range("B3").Font.Color = $hyperlinkColor

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Rows("2:2").Copy()
$wsZh.Rows("3:3").Insert()

$wsZh.Range("A2").Value = "5be07d11-a570-42c3-9722-9a26b9ee2579.md"
$wsZh.Range("G2").Value = "5be07d11-a570-42c3-9722-9a26b9ee2579.953d00e4d0d196bf98c720b67901336053b546ca.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-11-03 19:48:01"

$wsZh.Range("A3").Value = "ffff4be74221-375a-4846-b937-d78ae62f0c2f.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "5be07d11-a570-42c3-9722-9a26b9ee2579.953d00e4d0d196bf98c720b67901336053b546ca.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-11-03 19:48:01"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($repoBase + "5be07d11-a570-42c3-9722-9a26b9ee2579.md"), "", "", "5be07d11-a570-42c3-9722-9a26b9ee2579.md")
$wsZh.Range("A2").Font.Color = $hyperlinkColor

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($repoBase + "ffff4be74221-375a-4846-b937-d78ae62f0c2f.md"), "", "", "ffff4be74221-375a-4846-b937-d78ae62f0c2f.md")
$wsZh.Range("A3").Font.Color = $hyperlinkColor

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Rows("2:2").Copy()
$wsDe.Rows("3:3").Insert()

$wsDe.Range("A2").Value = "5be07d11-a570-42c3-9722-9a26b9ee2579.md"
$wsDe.Range("G2").Value = "5be07d11-a570-42c3-9722-9a26b9ee2579.953d00e4d0d196bf98c720b67901336053b546ca.de-de.xlf"
$wsDe.Range("H2").Value = "2016-11-03 19:48:14"

$wsDe.Range("A3").Value = "ffff4be74221-375a-4846-b937-d78ae62f0c2f.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "5be07d11-a570-42c3-9722-9a26b9ee2579.953d00e4d0d196bf98c720b67901336053b546ca.de-de.xlf"
$wsDe.Range("H3").Value = "2016-11-03 19:48:14"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($repoBase + "5be07d11-a570-42c3-9722-9a26b9ee2579.md"), "", "", "5be07d11-a570-42c3-9722-9a26b9ee2579.md")
$wsDe.Range("A2").Font.Color = $hyperlinkColor

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($repoBase + "ffff4be74221-375a-4846-b937-d78ae62f0c2f.md"), "", "", "ffff4be74221-375a-4846-b937-d78ae62f0c2f.md")
$wsDe.Range("A3").Font.Color = $hyperlinkColor

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
